$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 content (SC_003_01_ConsultaYDescargaClienteBancoIngles scenario) ---
$ws.Range("A2").Value = "SC_003_01_ConsultaYDescargaClienteBancoIngles"
$ws.Range("C2").Value = "Navegamos a la página principal y entramos en el home page.||Paso realizado correctamente."
$ws.Range("D2").Value = 'Seleccionar idioma Inglés y  la opción de "Idioma" .||Click realizado correctamente.'
$ws.Range("E2").Value = "Seleccionar el Rol Bank Client||Click al botón Bank Client."
$ws.Range("F2").Value = "Rellenamos los campos de datos, aceptamos condiciones y presionamos el botón Ok||Pasos realizados correctamente."
$ws.Range("G2").Value = "Seleccionamos el Año Fiscal  y presionar el botón Consultar||Pasos realizados correctamente."
$ws.Range("G2").VerticalAlignment = -4160
$ws.Range("G2").WrapText = $true
$ws.Range("H2").Value = "En la columna Type of request le damos click al icono de DOWNLOAD para descargar la constancia||Click realizado correctamente."
$ws.Range("H2").VerticalAlignment = -4160
$ws.Range("H2").WrapText = $true

# --- Row 3 content (SC_003_02_ConsultaYDescargaClienteBancoEspañol scenario) ---
$ws.Range("A3").Value = "SC_003_02_ConsultaYDescargaClienteBancoEspañol"
$ws.Range("C3").Value = "Navegamos a la página principal y entramos en el home page.||Paso realizado correctamente."
$ws.Range("E3").Value = "Seleccionar el Rol Cliente banco||Click al botón Cliente banco."
$ws.Range("D3").Value = 'Seleccionar idioma Español y la opción de "Idioma" .||Click realizado correctamente.'
$ws.Range("F3").Value = "Rellenamos los campos de datos, aceptamos condiciones y presionamos el botón Ok||Pasos realizados correctamente."
$ws.Range("G3").Value = "Seleccionamos el Año Fiscal  y presionar el botón Consultar||Pasos realizados correctamente."
$ws.Range("H3").Value = "En la columna Type of request le damos click al icono de DOWNLOAD para descargar la constancia||Click realizado correctamente."
$ws.Range("H3").VerticalAlignment = -4160
$ws.Range("H3").WrapText = $true

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 81.75
$ws.Rows.Item(3).RowHeight = 75

# --- Column widths (closest achievable values given COM ColumnWidth rounding) ---
$ws.Columns.Item(1).ColumnWidth = 44.833333333333336
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 25.833333333333332
$ws.Columns.Item(7).ColumnWidth = 27.5
$ws.Columns.Item(8).ColumnWidth = 27

# --- Selection / view ---
[void]$ws.Range("D3").Select()
